# fix hierarchy after testing and synthesis
#
# Adds a new "Solution 34: Hierarchical" row (row 38) to Sheet1, extending
# the results table that previously ended at row 37. Column layout mirrors
# the rows directly above it (35-37):
#   A: Solution name            F: Schedule/RTL phase
#   B: Throughput                G: Code-change note
#   C: Total area                H: =B*C            (Throughput*Area)
#   D: Pipelining note            I: =H/MIN($H$5:$H$110) (relative to smallest)
#   E: Unrolling note

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 38

$ws.Range("A$row").Value = "Solution 34: Hierarchical "
$ws.Range("B$row").Value = 64
$ws.Range("C$row").Value = 15653.38
$ws.Range("D$row").Value = "Main II=64"
$ws.Range("E$row").Value = "All the rest"
$ws.Range("F$row").Value = "RTL"
$ws.Range("G$row").Value = "Hierarchical design task 2"

$ws.Range("H$row").Formula = "=B$row*C$row"
$ws.Range("I$row").Formula = "=H$row/MIN(`$H`$5:`$H`$110)"

# Match the saved view state: scrolled up a bit with G36 as the active cell.
$ws.Range("G36").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1

Write-Output "Added row $row (Solution 34: Hierarchical) to $($ws.Name)"
